# Insert 4 new rows at row 804 (pushes old rows 804:847 down to 808:851)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("804:807").Insert(-4121)  # xlShiftDown

# New weekly data block (2021-11-16, serial 44516) for "Región del Maule"
# Columns: A Mercado ID, B Mercado, C Región, D Fecha, E Codreg, F Categoría ID,
#          G Categoría, H Variedad, I Calidad, J Volumen, K Precio min,
#          L Precio max, M Precio prom pond, N Unidad, O Origen, P Precio $/Kg,
#          Q Kg o Unidades, R Clasificación

$rows = @(
    @{ Row=804; H="Conconina(o)"; I="Primera"; J=500; K=4000; L=4000; M=4000; N="$/caja 10 unidades"; P=400; Q=10 },
    @{ Row=805; H="Escarola";     I="Primera"; J=800; K=4000; L=4000; M=4000; N="$/caja 15 unidades"; P=267; Q=15 },
    @{ Row=806; H="Española";    I="Primera"; J=500; K=4000; L=4000; M=4000; N="$/caja 18 unidades"; P=222; Q=18 },
    @{ Row=807; H="Marina";      I="Primera"; J=500; K=4000; L=4000; M=4000; N="$/caja 18 unidades"; P=222; Q=18 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 5
    $ws.Cells.Item($row, 2).Value = "Macroferia Regional de Talca"
    $ws.Cells.Item($row, 3).Value = "Maule"
    $ws.Cells.Item($row, 4).Value = 44516
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value = 7
    $ws.Cells.Item($row, 6).Value = 100112033
    $ws.Cells.Item($row, 7).Value = "Lechuga"
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = "Región del Maule"
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
